$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- Row 18: split/center-align the B18:C18 and G18:H18 pairs, then merge them --
$ws.Range("B18").Borders.Item(10).LineStyle = 0
$ws.Range("B18").HorizontalAlignment = -4108
$ws.Range("C18").Borders.Item(7).LineStyle = 0
$ws.Range("C18").HorizontalAlignment = -4108
$ws.Range("B18:C18").Merge()

$ws.Range("G18").Borders.Item(10).LineStyle = 0
$ws.Range("G18").HorizontalAlignment = -4108
$ws.Range("H18").Borders.Item(7).LineStyle = 0
$ws.Range("H18").HorizontalAlignment = -4108
$ws.Range("G18:H18").Merge()

# -- Text updates (also drives the sharedStrings reorder/relabel seen in the diff) --
$ws.Range("G18").Value = "Part 1"
$ws.Range("B18").Value = "Part 2"
$ws.Range("J3").Value = "init/Free"

# -- Cursor/selection position, matches the saved view state in the diff --
$ws.Range("L8").Select() | Out-Null
